$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells whose new values would otherwise be auto-parsed as numbers by Excel
# stay formatted as Text, matching the original inline-string cell content.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.657.68'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '3.137.27'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '571.15'
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").Value = '147.56'
$ws.Range("E6").Value = '  -1.40%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").Value = '3.137.92'
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("D9").Value = '0.521'
$ws.Range("E9").Value = '  -0.67%  '
$ws.Range("D10").Value = '0.156'
$ws.Range("E10").Value = '  -3.73%  '
$ws.Range("D11").Value = '6.03'
$ws.Range("E11").Value = '  -2.25%  '
$ws.Range("D12").Value = '0.492'
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").Value = '  +2.64%  '
$ws.Range("D14").Value = '36.63'
$ws.Range("E14").Value = '  -1.25%  '
$ws.Range("D15").Value = '3.650.09'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '64.792.72'
$ws.Range("E16").Value = '  -0.44%  '
$ws.Range("D17").Value = '3.140.36'
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("D18").Value = '7.02'
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("D20").Value = '497.37'
$ws.Range("E20").Value = '  -1.77%  '
$ws.Range("D21").Value = '14.68'
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").Value = '0.706'
$ws.Range("E22").Value = '  -2.09%  '
$ws.Range("D23").Value = '15.04'
$ws.Range("E23").Value = '  -3.84%  '
$ws.Range("D24").Value = '7.61'
$ws.Range("E24").Value = '  -2.09%  '
$ws.Range("D25").Value = '83.39'
$ws.Range("E25").Value = '  -1.80%  '
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").Value = '2.86'
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '8.76'
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("E29").Value = '  +0.46%  '
$ws.Range("E30").Value = '  +3.01%  '
$ws.Range("D31").Value = '27.26'
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("E33").Value = '  -0.33%  '
$ws.Range("D34").Value = '6.11'
$ws.Range("E34").Value = '  +1.48%  '
$ws.Range("D35").Value = '6.38'
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("D36").Value = '54.38'
$ws.Range("E36").Value = '  -2.46%  '
$ws.Range("D37").Value = '0.0888'
$ws.Range("E37").Value = '  +4.49%  '
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("D39").Value = '0.0410'
$ws.Range("E39").Value = '  -3.22%  '
$ws.Range("D40").Value = '2.93'
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("D41").Value = '8.56'
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").Value = '3.009.51'
$ws.Range("E42").Value = '  -3.42%  '
$ws.Range("E43").Value = '  -4.76%  '
$ws.Range("D44").Value = '0.279'
$ws.Range("E44").Value = '  -2.63%  '
$ws.Range("D45").Value = '2.39'
$ws.Range("E45").Value = '  -0.58%  '
$ws.Range("D46").Value = '27.91'
$ws.Range("E46").Value = '  -3.88%  '
$ws.Range("D47").Value = '0.0₃0568'
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("E49").Value = '  -1.94%  '
$ws.Range("D50").Value = '2.21'
$ws.Range("E50").Value = '  -2.09%  '
$ws.Range("D51").Value = '117.51'
$ws.Range("E51").Value = '  -0.78%  '
